$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset all existing formatting/content in the used area to a clean slate,
# so newly written cells start from the true default style (no "s" attr).
$ws.Range("A1:G6").EntireRow.Delete()
$ws.Range("A1:G6").EntireColumn.Delete()

# ---- Header row (row 1) ----
$ws.Range("A1").Value2 = "Wave"
$ws.Range("B1").Value2 = "MonsterClassPath"
$ws.Range("C1").Value2 = "Num"
$ws.Range("D1").Value2 = "Damage"
$ws.Range("E1").Value2 = "MaxHP"
$ws.Range("F1").Value2 = "Speed"
$ws.Range("G1").Value2 = "Size"
$ws.Range("H1").Value2 = "DropMoney"

# ---- Row 2 ----
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "(""/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C"")"
$ws.Range("C2").Value2 = -9
$ws.Range("D2").Value2 = -50
$ws.Range("E2").Value2 = -10
$ws.Range("F2").Value2 = -200
$ws.Range("G2").Value2 = -1
$ws.Range("H2").Value2 = -10

# ---- Row 3 ----
$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "(""/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C"")"
$ws.Range("C3").Value2 = -20
$ws.Range("D3").Value2 = -10
$ws.Range("E3").Value2 = -50
$ws.Range("F3").Value2 = -400
$ws.Range("G3").Value2 = -1
$ws.Range("H3").Value2 = -100

# ---- Row 4 ----
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "(""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C"")"
$ws.Range("C4").Value2 = "(1, 2, 3)"
$ws.Range("D4").Value2 = "(50, 20, 50)"
$ws.Range("E4").Value2 = "(100, 200, 300)"
$ws.Range("F4").Value2 = "(200, 300, 400)"
$ws.Range("G4").Value2 = "(1, 1, 1)"
$ws.Range("H4").Value2 = "(1, 1, 1)"

# ---- Row 5 ----
$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = "(""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C"")"
$ws.Range("C5").Value2 = "(1, 2, 3)"
$ws.Range("D5").Value2 = "(50, 20, 50)"
$ws.Range("E5").Value2 = "(100, 200, 300)"
$ws.Range("F5").Value2 = "(200, 300, 400)"
$ws.Range("G5").Value2 = "(1, 1, 1)"
$ws.Range("H5").Value2 = "(1, 1, 1)"

# ---- Row 6 ----
$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value2 = "(""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/SuperMinion/BSuperMinion_BP.BSuperMinion_BP_C"", ""/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C"")"
$ws.Range("C6").Value2 = "(1, 2, 3)"
$ws.Range("D6").Value2 = "(50, 20, 50)"
$ws.Range("E6").Value2 = "(100, 200, 300)"
$ws.Range("F6").Value2 = "(200, 300, 400)"
$ws.Range("G6").Value2 = "(1, 1, 1)"
$ws.Range("H6").Value2 = "(1, 1, 1)"

# ---- Formatting ----
# Column A + header row: centered, no wrap (style "1" in the target)
$ws.Range("A1:A6").HorizontalAlignment = -4108
$ws.Range("A1:A6").WrapText = $false
$ws.Range("A1:H1").HorizontalAlignment = -4108
$ws.Range("A1:H1").WrapText = $false

# Boss rows (2-3) numeric columns: custom number format, default alignment
$ws.Range("C2:H3").NumberFormat = "0_);\(0\)"

# Column widths (engine stores width as whole pixels at 7px/char + 5px
# padding, so request the char-width values whose rounded-pixel result is
# closest to the target stored widths of 6.125 and 72)
$ws.Columns.Item(1).ColumnWidth = 5.428571428571429
$ws.Columns.Item(2).ColumnWidth = 71.28571428571429

# Row heights back to auto (drop explicit/custom row height)
$ws.Range("A1:H6").EntireRow.AutoFit()

# Selection / view
$ws.Range("B2").Select()
